$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WPP")

$ws.Range("D8").Value = 20613600
$ws.Range("E8").Value = 18767600
$ws.Range("F8").Value = 15958500
$ws.Range("G8").Value = 15037300
$ws.Range("H8").Value = 14372700
$ws.Range("I8").Value = 13529700
$ws.Range("J8").Value = 13071500
$ws.Range("D9").Value = 16472100
$ws.Range("E9").Value = 14801400
$ws.Range("F9").Value = 12663700
$ws.Range("G9").Value = 1909600
$ws.Range("H9").Value = 1230400
$ws.Range("I9").Value = 1119500
$ws.Range("J9").Value = 1021700
$ws.Range("D10").Value = 4141400
$ws.Range("E10").Value = 3966100
$ws.Range("F10").Value = 3294800
$ws.Range("G10").Value = 13127600
$ws.Range("H10").Value = 13142400
$ws.Range("I10").Value = 12410200
$ws.Range("J10").Value = 12049900
$ws.Range("E14").Value = -119900
$ws.Range("F14").Value = 114900
$ws.Range("G14").Value = 244700
$ws.Range("H14").Value = 47300
$ws.Range("I14").Value = 187600
$ws.Range("E15").Value = 550700
$ws.Range("F15").Value = 474500
$ws.Range("G15").Value = 483600
$ws.Range("H15").Value = 532200
$ws.Range("I15").Value = 508400
$ws.Range("J15").Value = 490900
$ws.Range("D17").Value = 18124700
$ws.Range("E17").Value = 16076700
$ws.Range("F17").Value = 13829900
$ws.Range("G17").Value = 13081200
$ws.Range("H17").Value = 12547200
$ws.Range("I17").Value = 11911000
$ws.Range("J17").Value = 11516500
$ws.Range("D18").Value = 2488900
$ws.Range("E18").Value = 2690900
$ws.Range("F18").Value = 2128600
$ws.Range("G18").Value = 1956100
$ws.Range("H18").Value = 1825500
$ws.Range("I18").Value = 1618800
$ws.Range("J18").Value = 1555000
$ws.Range("D20").Value = 600900
$ws.Range("E20").Value = 94600
$ws.Range("F20").Value = 97700
$ws.Range("G20").Value = 267400
$ws.Range("H20").Value = 197000
$ws.Range("I20").Value = 140300
$ws.Range("J20").Value = 88400
$ws.Range("I21").Value = 2282200
$ws.Range("J21").Value = "NA"
$ws.Range("D22").Value = 338600
$ws.Range("E22").Value = 319700
$ws.Range("F22").Value = 279500
$ws.Range("G22").Value = 329700
$ws.Range("H22").Value = 332300
$ws.Range("I22").Value = 334900
$ws.Range("J22").Value = 328200
$ws.Range("D23").Value = 2751200
$ws.Range("E23").Value = 2465800
$ws.Range("F23").Value = 1946800
$ws.Range("G23").Value = 1893700
$ws.Range("H23").Value = 1690100
$ws.Range("I23").Value = 1424200
$ws.Range("J23").Value = 1315300
$ws.Range("D24").Value = 525600
$ws.Range("E24").Value = 507200
$ws.Range("F24").Value = 322800
$ws.Range("G24").Value = 391800
$ws.Range("H24").Value = 370000
$ws.Range("I24").Value = 257200
$ws.Range("J24").Value = 119900
$ws.Range("D26").Value = 2225500
$ws.Range("E26").Value = 1958600
$ws.Range("F26").Value = 1624000
$ws.Range("G26").Value = 1501900
$ws.Range("H26").Value = 1320100
$ws.Range("I26").Value = 1167000
$ws.Range("J26").Value = 1195400
$ws.Range("D27").Value = 2100700
$ws.Range("E27").Value = 1826200
$ws.Range("F27").Value = 1513300
$ws.Range("G27").Value = 1405000
$ws.Range("H27").Value = 1221500
$ws.Range("I27").Value = 1073100
$ws.Range("J27").Value = 1095800
$ws.Range("D29").Value = 268700
$ws.Range("D32").Value = -600900
$ws.Range("E32").Value = -94600
$ws.Range("F32").Value = -97700
$ws.Range("G32").Value = -267400
$ws.Range("H32").Value = -197000
$ws.Range("I32").Value = -140300
$ws.Range("J32").Value = -88400
$ws.Range("D33").Value = 2369400
$ws.Range("E33").Value = 1826200
$ws.Range("F33").Value = 1513300
$ws.Range("G33").Value = 1405000
$ws.Range("H33").Value = 1221500
$ws.Range("I33").Value = 1073100
$ws.Range("J33").Value = 1095800
$ws.Range("D35").Value = 2369400
$ws.Range("E35").Value = 1826200
$ws.Range("F35").Value = 1513300
$ws.Range("G35").Value = 1405000
$ws.Range("H35").Value = 1221500
$ws.Range("I35").Value = 1073100
$ws.Range("J35").Value = 1095800
$ws.Range("D41").Value = 3119100
$ws.Range("E41").Value = 2942800
$ws.Range("F41").Value = 2905700
$ws.Range("G41").Value = 2565600
$ws.Range("H41").Value = 2737900
$ws.Range("I41").Value = 2245200
$ws.Range("J41").Value = 2391500
$ws.Range("E42").Value = 235700
$ws.Range("F42").Value = 201600
$ws.Range("G42").Value = 711800
$ws.Range("H42").Value = 159800
$ws.Range("I42").Value = 292000
$ws.Range("J42").Value = 147500
$ws.Range("D43").Value = 15736500
$ws.Range("E43").Value = 16018200
$ws.Range("F43").Value = 13596600
$ws.Range("G43").Value = 12315400
$ws.Range("H43").Value = 11628100
$ws.Range("I43").Value = 8998400
$ws.Range("J43").Value = 9083500
$ws.Range("D44").Value = 553400
$ws.Range("E44").Value = 522200
$ws.Range("F44").Value = 429100
$ws.Range("G44").Value = 426900
$ws.Range("H44").Value = 397200
$ws.Range("I44").Value = 454200
$ws.Range("J44").Value = 435500
$ws.Range("D45").Value = 390400
$ws.Range("E45").Value = 423500
$ws.Range("F45").Value = 312500
$ws.Range("G45").Value = 304600
$ws.Range("H45").Value = 403000
$ws.Range("I45").Value = 2911500
$ws.Range("J45").Value = 2666000
$ws.Range("D46").Value = 19769200
$ws.Range("E46").Value = 20142500
$ws.Range("F46").Value = 17445700
$ws.Range("G46").Value = 16324200
$ws.Range("H46").Value = 15325900
$ws.Range("I46").Value = 14901400
$ws.Range("J46").Value = 14724000
$ws.Range("D47").Value = 2893900
$ws.Range("E47").Value = 3103900
$ws.Range("F47").Value = 2500800
$ws.Range("G47").Value = 1864000
$ws.Range("H47").Value = 1593700
$ws.Range("I47").Value = 1483800
$ws.Range("J47").Value = 1408300
$ws.Range("D48").Value = 2555100
$ws.Range("E48").Value = 1263500
$ws.Range("F48").Value = 1040400
$ws.Range("G48").Value = 1007600
$ws.Range("H48").Value = 1008600
$ws.Range("I48").Value = 1002100
$ws.Range("J48").Value = 949900
$ws.Range("D49").Value = 39054400
$ws.Range("E49").Value = 20127600
$ws.Range("F49").Value = 16155200
$ws.Range("G49").Value = 15193000
$ws.Range("H49").Value = 14530800
$ws.Range("I49").Value = 14718600
$ws.Range("J49").Value = 14726600
$ws.Range("D52").Value = 438900
$ws.Range("E52").Value = 450400
$ws.Range("F52").Value = 355800
$ws.Range("G52").Value = 506500
$ws.Range("H52").Value = 362500
$ws.Range("I52").Value = 342300
$ws.Range("J52").Value = 401100
$ws.Range("D54").Value = 43906700
$ws.Range("E54").Value = 45087800
$ws.Range("F54").Value = 37497900
$ws.Range("G54").Value = 34895200
$ws.Range("H54").Value = 32614800
$ws.Range("I54").Value = 32448100
$ws.Range("J54").Value = 32209800
$ws.Range("D57").Value = 13139200
$ws.Range("E57").Value = 13807200
$ws.Range("F57").Value = 11300900
$ws.Range("G57").Value = 10321500
$ws.Range("H57").Value = 9390900
$ws.Range("I57").Value = 9470500
$ws.Range("J57").Value = 9638200
$ws.Range("D58").Value = 814000
$ws.Range("E58").Value = 1307600
$ws.Range("F58").Value = 1215600
$ws.Range("G58").Value = 852000
$ws.Range("H58").Value = 2455800
$ws.Range("I58").Value = 1416400
$ws.Range("J58").Value = 676200
$ws.Range("D59").Value = 6282500
$ws.Range("E59").Value = 6752300
$ws.Range("F59").Value = 6024900
$ws.Range("G59").Value = 5255300
$ws.Range("H59").Value = 4735800
$ws.Range("I59").Value = 4890900
$ws.Range("J59").Value = 5073000
$ws.Range("D60").Value = 20235700
$ws.Range("E60").Value = 21867000
$ws.Range("F60").Value = 18541400
$ws.Range("G60").Value = 16428800
$ws.Range("H60").Value = 15354600
$ws.Range("I60").Value = 15777700
$ws.Range("J60").Value = 15387300
$ws.Range("D61").Value = 8152500
$ws.Range("E61").Value = 7258400
$ws.Range("F61").Value = 6079600
$ws.Range("G61").Value = 5393200
$ws.Range("H61").Value = 4592000
$ws.Range("I61").Value = 4800600
$ws.Range("J61").Value = 5077700
$ws.Range("D62").Value = 2532700
$ws.Range("E62").Value = 3222400
$ws.Range("F62").Value = 2421700
$ws.Range("G62").Value = 2864700
$ws.Range("H62").Value = 2434000
$ws.Range("I62").Value = 2660500
$ws.Range("J62").Value = 2752500
$ws.Range("D66").Value = 31532300
$ws.Range("E66").Value = 32925700
$ws.Range("F66").Value = 27536300
$ws.Range("G66").Value = 25112800
$ws.Range("H66").Value = 22719300
$ws.Range("I66").Value = 23564400
$ws.Range("J66").Value = 23518800
$ws.Range("D72").Value = 12654500
$ws.Range("E72").Value = 10645100
$ws.Range("F72").Value = 9751000
$ws.Range("G72").Value = 8951700
$ws.Range("H72").Value = 8845300
$ws.Range("I72").Value = 7759900
$ws.Range("J72").Value = 7278600
$ws.Range("D76").Value = 12374400
$ws.Range("E76").Value = 12162000
$ws.Range("F76").Value = 9961500
$ws.Range("G76").Value = 9782500
$ws.Range("H76").Value = 9895500
$ws.Range("I76").Value = 8883700
$ws.Range("J76").Value = 8691000
$ws.Range("D81").Value = 2369400
$ws.Range("E81").Value = 1826200
$ws.Range("F81").Value = 1513300
$ws.Range("G81").Value = 1405000
$ws.Range("H81").Value = 1221500
$ws.Range("I81").Value = 1073100
$ws.Range("J81").Value = 1095800
$ws.Range("D83").Value = 602700
$ws.Range("E83").Value = 558000
$ws.Range("F83").Value = 480600
$ws.Range("G83").Value = 490900
$ws.Range("H83").Value = 540600
$ws.Range("I83").Value = 517300
$ws.Range("J83").Value = "NA"
$ws.Range("D89").Value = 1836600
$ws.Range("E89").Value = 2313600
$ws.Range("F89").Value = 1773700
$ws.Range("G89").Value = 2222200
$ws.Range("H89").Value = 1792400
$ws.Range("I89").Value = 1184700
$ws.Range("J89").Value = 867600
$ws.Range("D91").Value = -376800
$ws.Range("E91").Value = -328800
$ws.Range("F91").Value = -274300
$ws.Range("G91").Value = -232000
$ws.Range("H91").Value = -313900
$ws.Range("I91").Value = -378600
$ws.Range("J91").Value = -281900
$ws.Range("D94").Value = -651800
$ws.Range("E94").Value = -1195000
$ws.Range("F94").Value = -1177100
$ws.Range("G94").Value = -909900
$ws.Range("H94").Value = -624200
$ws.Range("I94").Value = -1008400
$ws.Range("J94").Value = "NA"
$ws.Range("D96").Value = -980200
$ws.Range("E96").Value = -804100
$ws.Range("F96").Value = -711900
$ws.Range("G96").Value = -600000
$ws.Range("H96").Value = -518200
$ws.Range("I96").Value = -399900
$ws.Range("J96").Value = -284900
$ws.Range("D100").Value = -1024700
$ws.Range("E100").Value = -1556700
$ws.Range("F100").Value = -918200
$ws.Range("G100").Value = -745300
$ws.Range("H100").Value = -270000
$ws.Range("I100").Value = -110600
$ws.Range("J100").Value = "NA"
$ws.Range("D101").Value = -35500
$ws.Range("E101").Value = 380700
$ws.Range("F101").Value = -71000
$ws.Range("G101").Value = -91700
$ws.Range("H101").Value = -214800
$ws.Range("I101").Value = -155600
$ws.Range("J101").Value = "NA"
$ws.Range("D102").Value = 124700
$ws.Range("E102").Value = -57400
$ws.Range("F102").Value = -392600
$ws.Range("G102").Value = 475300
$ws.Range("H102").Value = 683300
$ws.Range("I102").Value = -89900
$ws.Range("J102").Value = -367300
